$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Metadata")
$ws1.Range("B8").Value = "2025-07-11T12:29:53+00:00"
$ws1.Range("B11").Value = "FRANCE"
